$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 3 (shifting existing data down),
# completing the import of a new defect/description row for ticket 8458.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row with the new defect description.
$ws.Range("A3").Value = "משהו חדש שלי"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 8458

# Leave the selection where the author left off editing.
$ws.Range("H5").Select()
